# Add new worksheet "U chart" with sample_size data (for SVM pattern detection in control chart app)

$wb = $excel.ActiveWorkbook

# --- Add the new "U chart" worksheet as the last sheet ---
$wsCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($wsCount)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "U chart"

# --- Populate header row ---
$newSheet.Range("A1").Value = "Sample"
$newSheet.Range("B1").Value = 1
$newSheet.Range("C1").Value = "sample_size"

# --- Populate data rows ---
$data = @(
    @(1, 52, 100),
    @(2, 48, 120),
    @(3, 56, 150),
    @(4, 25, 80),
    @(5, 39, 130),
    @(6, 39, 125),
    @(7, 54, 140),
    @(8, 70, 90),
    @(9, 41, 100),
    @(10, 43, 110),
    @(11, 47, 115),
    @(12, 52, 120),
    @(13, 44, 90),
    @(14, 47, 85),
    @(15, 50, 100),
    @(16, 40, 120),
    @(17, 47, 115),
    @(18, 46, 120),
    @(19, 44, 120),
    @(20, 50, 110)
)

$row = 2
foreach ($r in $data) {
    $newSheet.Cells.Item($row, 1).Value = $r[0]
    $newSheet.Cells.Item($row, 2).Value = $r[1]
    $newSheet.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# --- Set selection on the new sheet ---
$newSheet.Range("D2").Select()

# --- Update selections on other sheets ---
$ws1 = $wb.Worksheets.Item("X bar and R chart")
$ws1.Range("I16").Select()

$ws2 = $wb.Worksheets.Item("C chart")
$ws2.Range("E27").Select()

# --- Make "U chart" the active sheet/tab ---
$newSheet.Activate()
